$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "1 on 09/09"
$ws.Range("K2").Value = "1 on 09/10"
$ws.Range("F4").Value = "1 on 09/10"
$ws.Range("I6").Value = "1 on 09/09"
$ws.Range("G8").Value = "1 on 09/10"
$ws.Range("G10").Value = "1 on 09/09"
$ws.Range("E12").Value = "1 on 09/10"

# Column K now has content for the first time, so Excel auto-sizes it to
# fit "1 on 09/10" (best-fit width of 10 characters).
$ws.Columns.Item(11).AutoFit() | Out-Null
$ws.Columns.Item(11).ColumnWidth = 9.17

$ws.Range("E12").Select()
